# Commit: "atualizacao do material sobre ppo"
#
# Semantic edit on Sheet1 (the "21-Maio" / final-evaluation block):
#  - Row 30 (21-Maio) incorrectly repeated "Apresentação final do projeto"
#    in C:F; it should carry the final-evaluation text that used to live
#    two rows further down (2-Junho/4-Junho), which are being removed.
#  - Rows 31-32 (26-Maio/28-Maio) stay "Não teremos aula" (unchanged).
#  - Rows 33-34 (2-Junho/4-Junho) are removed from the schedule: their
#    dates and contents are cleared out (course now ends earlier).
#  - The two now-permanently-blank rows that used to sit right after
#    (old rows 35-36) are deleted so everything below shifts up by two,
#    closing the gap left by the removed 2-Junho/4-Junho rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Fix row 30 (21-Maio): replace the wrongly duplicated text with the
#    final-assessment wording that used to sit in row 33.
$ws.Range("C30").Value = "Avaliação final"
$ws.Range("D30").Value = "Avaliação"
$ws.Range("E30").Value = "Avaliação Final"
$ws.Range("F30").Value = "Avaliação Final"

# 2) Clear out the old "2-Junho" (row 33) and "4-Junho" (row 34) rows --
#    date column A and content columns C:F. Column B stays empty as before.
$ws.Range("A33:A34,C33:F34").ClearContents()

# Re-assert the "plain / no special fill" look on the cleared cells so the
# formatting matches the rest of the blank cells on the sheet (no leftover
# weekday-row shading or the old "Bad" red styling).
foreach ($addr in @("A33","A34")) {
    $cell = $ws.Range($addr)
    $cell.Interior.ColorIndex = -4142
    $cell.WrapText = $true
    $cell.HorizontalAlignment = -4142
}
foreach ($addr in @("C33","D33","E33","F33","C34","D34","E34","F34")) {
    $cell = $ws.Range($addr)
    $cell.WrapText = $true
    $cell.HorizontalAlignment = -4142
}

# 3) Close the gap: the two rows right after (old 35 & 36) are already
#    completely empty placeholders -- deleting them shifts the remaining
#    tail (old 37-39) up by two, becoming the new 35-37.
$ws.Range("A35:A36").EntireRow.Delete()
